# Auto-generated edit script applying the Sargatanas_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N19").Value = -9019
$ws.Range("H19").Value = 5893.6
$ws.Range("L19").Value = 8669
$ws.Range("J19").Value = 8669
$ws.Range("L53").Value = 4816.143
$ws.Range("J53").Value = 4816.143
$ws.Range("N53").Value = -6090.143
$ws.Range("H53").Value = 4138.3125
$ws.Range("H96").Value = 997.375
$ws.Range("K96").Value = 2492.4999
$ws.Range("M96").Value = -1119.4999
$ws.Range("I96").Value = 830.8333
$ws.Range("K98").Value = 4312.0713
$ws.Range("M98").Value = -2814.0713
$ws.Range("I98").Value = 4312.0713
$ws.Range("H98").Value = 4091.6667
$ws.Range("I113").Value = 22224246
$ws.Range("M113").Value = -22220992
$ws.Range("H113").Value = 38200484
$ws.Range("L113").Value = 45462410
$ws.Range("K113").Value = 22224246
$ws.Range("N113").Value = -45468918
$ws.Range("J113").Value = 45462410
$ws.Range("K115").Value = 1632.6
$ws.Range("I115").Value = 544.2
$ws.Range("H115").Value = 544.2
$ws.Range("M115").Value = -65.60000000000014
$ws.Range("K122").Value = 12936.2139
$ws.Range("H122").Value = 4091.6667
$ws.Range("M122").Value = -10486.2139
$ws.Range("I122").Value = 4312.0713
$ws.Range("M131").Value = -3040.875
$ws.Range("J131").Value = 4499.6665
$ws.Range("L131").Value = 13498.9995
$ws.Range("N131").Value = -23578.9995
$ws.Range("K131").Value = 8080.875
$ws.Range("H131").Value = 3186.182
$ws.Range("I131").Value = 2693.625
$ws.Range("H137").Value = 2508.1
$ws.Range("L137").Value = 9530.3079
$ws.Range("J137").Value = 3176.7693
$ws.Range("N137").Value = -14630.3079
$ws.Range("I137").Value = 2186.1482
$ws.Range("K137").Value = 6558.444600000001
$ws.Range("M137").Value = -4008.444600000001
$ws.Range("H138").Value = 4265.3022
$ws.Range("N138").Value = -38310.6875
$ws.Range("L138").Value = 28030.6875
$ws.Range("J138").Value = 9343.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K2").Value = 2398.6
$ws.Range("H2").Value = 4969.3
$ws.Range("I2").Value = 2398.6
$ws.Range("M2").Value = -2285.6
$ws.Range("K61").Value = 1794.1892
$ws.Range("H61").Value = 4616.66
$ws.Range("N61").Value = -13073.846
$ws.Range("J61").Value = 12649.846
$ws.Range("L61").Value = 12649.846
$ws.Range("M61").Value = -1582.1892
$ws.Range("I61").Value = 1794.1892
$ws.Range("J63").Value = 2343.75
$ws.Range("N63").Value = -3715.75
$ws.Range("H63").Value = 2356.625
$ws.Range("L63").Value = 2343.75
$ws.Range("N66").Value = -18582.75
$ws.Range("H66").Value = 2356.625
$ws.Range("L66").Value = 11718.75
$ws.Range("J66").Value = 2343.75
$ws.Range("M116").Value = -104.5999999999999
$ws.Range("I116").Value = 2398.6
$ws.Range("K116").Value = 2398.6
$ws.Range("H116").Value = 4969.3
$ws.Range("K122").Value = 208510.992
$ws.Range("J122").Value = 8000
$ws.Range("N122").Value = -28900
$ws.Range("H122").Value = 31063.875
$ws.Range("M122").Value = -206060.992
$ws.Range("I122").Value = 69503.664
$ws.Range("L122").Value = 24000
$ws.Range("J136").Value = 12649.846
$ws.Range("L136").Value = 37949.538
$ws.Range("N136").Value = -43049.538
$ws.Range("K136").Value = 5382.5676
$ws.Range("I136").Value = 1794.1892
$ws.Range("H136").Value = 4616.66
$ws.Range("M136").Value = -2832.5676
$ws.Range("N139").Value = -90255
$ws.Range("L139").Value = 79975
$ws.Range("H139").Value = 79975
$ws.Range("J139").Value = 79975

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M3").Value = -2284.6
$ws.Range("H3").Value = 4969.3
$ws.Range("I3").Value = 2398.6
$ws.Range("K3").Value = 2398.6
$ws.Range("N105").Value = -7607.4287
$ws.Range("H105").Value = 52552.16
$ws.Range("I105").Value = 66680.125
$ws.Range("M105").Value = -64933.125
$ws.Range("K105").Value = 66680.125
$ws.Range("J105").Value = 4113.4287
$ws.Range("L105").Value = 4113.4287
$ws.Range("N107").Value = -10973.3335
$ws.Range("L107").Value = 7133.3335
$ws.Range("J107").Value = 7133.3335
$ws.Range("H107").Value = 51139004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M31").Value = -3246.1
$ws.Range("H31").Value = 8430
$ws.Range("I31").Value = 3541.1
$ws.Range("K31").Value = 3541.1
$ws.Range("M34").Value = -3339.1
$ws.Range("H34").Value = 8430
$ws.Range("I34").Value = 3541.1
$ws.Range("K34").Value = 3541.1
$ws.Range("K122").Value = 3480
$ws.Range("J122").Value = 2296.625
$ws.Range("N122").Value = -11789.875
$ws.Range("H122").Value = 1728.3125
$ws.Range("M122").Value = -1030
$ws.Range("I122").Value = 1160
$ws.Range("L122").Value = 6889.875
$ws.Range("K132").Value = 5866.9998
$ws.Range("M132").Value = -3336.9998
$ws.Range("N132").Value = -35108
$ws.Range("L132").Value = 30048
$ws.Range("J132").Value = 10016
$ws.Range("I132").Value = 1955.6666
$ws.Range("H132").Value = 6791.8667
$ws.Range("I134").Value = 1643.6666
$ws.Range("K134").Value = 4930.9998
$ws.Range("M134").Value = -2395.9998
$ws.Range("H134").Value = 6788.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K14").Value = 50001909
$ws.Range("M14").Value = -50001736
$ws.Range("H14").Value = 16667303
$ws.Range("I14").Value = 16667303
$ws.Range("H34").Value = 5972.5713
$ws.Range("L34").Value = 19290.693
$ws.Range("J34").Value = 6430.231
$ws.Range("N34").Value = -19458.693
$ws.Range("H87").Value = 888.3333
$ws.Range("K87").Value = 2664.9999
$ws.Range("I87").Value = 888.3333
$ws.Range("M87").Value = -1416.9999
$ws.Range("I90").Value = 888.3333
$ws.Range("K90").Value = 7994.9997
$ws.Range("H90").Value = 888.3333
$ws.Range("M90").Value = -1754.9997
$ws.Range("H137").Value = 144629.58
$ws.Range("I137").Value = 143846.72
$ws.Range("K137").Value = 431540.16
$ws.Range("M137").Value = -426440.16

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N80").Value = -6719.375
$ws.Range("J80").Value = 4723.375
$ws.Range("H80").Value = 4435.091
$ws.Range("I80").Value = 3666.3333
$ws.Range("L80").Value = 4723.375
$ws.Range("M80").Value = -2668.3333
$ws.Range("K80").Value = 3666.3333
$ws.Range("N83").Value = -33600.875
$ws.Range("H83").Value = 4435.091
$ws.Range("K83").Value = 18331.6665
$ws.Range("M83").Value = -13339.6665
$ws.Range("J83").Value = 4723.375
$ws.Range("I83").Value = 3666.3333
$ws.Range("L83").Value = 23616.875
$ws.Range("M102").Value = 60.71419999999989
$ws.Range("K102").Value = 1561.2858
$ws.Range("H102").Value = 1672.2941
$ws.Range("I102").Value = 1561.2858
$ws.Range("K122").Value = 18114057
$ws.Range("J122").Value = 4388.143
$ws.Range("N122").Value = -18064.429
$ws.Range("H122").Value = 3815102.2
$ws.Range("M122").Value = -18111607
$ws.Range("I122").Value = 6038019
$ws.Range("L122").Value = 13164.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N7").Value = -9391.5
$ws.Range("H7").Value = 4423.591
$ws.Range("K7").Value = 2644.625
$ws.Range("L7").Value = 9167.5
$ws.Range("I7").Value = 2644.625
$ws.Range("M7").Value = -2532.625
$ws.Range("J7").Value = 9167.5
$ws.Range("J40").Value = 5667.0835
$ws.Range("N40").Value = -5939.0835
$ws.Range("H40").Value = 5307.4
$ws.Range("L40").Value = 5667.0835
$ws.Range("I40").Value = 4767.875
$ws.Range("K40").Value = 4767.875
$ws.Range("M40").Value = -4631.875
$ws.Range("K126").Value = 7933.875
$ws.Range("H126").Value = 4423.591
$ws.Range("N126").Value = -32442.5
$ws.Range("M126").Value = -5463.875
$ws.Range("L126").Value = 27502.5
$ws.Range("I126").Value = 2644.625
$ws.Range("J126").Value = 9167.5
$ws.Range("J136").Value = 13149.95
$ws.Range("L136").Value = 39449.85000000001
$ws.Range("N136").Value = -44549.85000000001
$ws.Range("K136").Value = 3024.2307
$ws.Range("I136").Value = 1008.0769
$ws.Range("H136").Value = 8366.788
$ws.Range("M136").Value = -474.2307000000001
$ws.Range("H138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 728
$ws.Range("I100").Value = 542.5294
$ws.Range("K100").Value = 1085.0588
$ws.Range("M100").Value = -544.0588
$ws.Range("K122").Value = 805307.3999999999
$ws.Range("H122").Value = 169324.8
$ws.Range("M122").Value = -802857.3999999999
$ws.Range("I122").Value = 268435.8
$ws.Range("K132").Value = 48398529
$ws.Range("M132").Value = -48395999
$ws.Range("N132").Value = -93924.5
$ws.Range("L132").Value = 88864.5
$ws.Range("J132").Value = 29621.5
$ws.Range("I132").Value = 16132843
$ws.Range("H132").Value = 12205228
$ws.Range("J136").Value = 49059.375
$ws.Range("L136").Value = 147178.125
$ws.Range("N136").Value = -152278.125
$ws.Range("K136").Value = 88237332
$ws.Range("I136").Value = 29412444
$ws.Range("H136").Value = 17262078
$ws.Range("M136").Value = -88234782
